$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("F").Delete()
foreach ($hl in $ws.Hyperlinks) { $hl.Delete() }
$rng = $ws.Range("H2")
$rng.Hyperlinks.Add($rng, "mailto:olo@gmail.com", "", "", "olo@gmail.com")
$rng.Font.Name = "Calibri"
$rng.Font.Name = "Arial"
$rng.Font.Color = 16711680
$rng.Font.Underline = -4142
